$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edges")

# Update the edge count cell (graph is now bidirectional, so edge count doubled)
$ws.Range("A1").Value = 18

# Append the reverse direction of each existing edge (rows 2-10) as new rows 11-19
$newEdges = @(
    @(5,1),
    @(6,1),
    @(7,1),
    @(5,2),
    @(8,2),
    @(5,3),
    @(5,3),
    @(8,3),
    @(8,4)
)

$row = 11
foreach ($edge in $newEdges) {
    $ws.Cells.Item($row, 1).Value = $edge[0]
    $ws.Cells.Item($row, 2).Value = $edge[1]
    $row++
}

$ws.Activate()
$ws.Range("B20").Select()
